# edit.ps1
# Applies the weekly CompStat data refresh described by the commit
# "New crime data collected": updated report-date strings in the title
# block plus refreshed Murder..Shooting Vic. crime-count/percent-change
# figures for rows 15-27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PlainValue($cell, $val) {
    $ws.Range($cell).Value = $val
}

function Set-TextValue($cell, $val) {
    # Some "count" cells in this report are stored as literal text
    # (e.g. "0" or "***.*") rather than as numbers. Force the Text
    # number format while writing so the numeric-looking string is
    # not auto-coerced back into a number, then restore General so the
    # cell keeps its original (right-aligned, unformatted) look.
    $ws.Range($cell).NumberFormat = "@"
    $ws.Range($cell).Value = $val
    $ws.Range($cell).NumberFormat = "general"
}

function Set-NumValue($cell, $val, $fmt) {
    # Cell switches from the literal-text representation back to a
    # real number; (re)apply the column's numeric format explicitly.
    $ws.Range($cell).NumberFormat = $fmt
    $ws.Range($cell).Value = $val
}

# --- Header / report-period strings --------------------------------
$ws.Range("A8").Value = "Volume 30   Number  10"
$ws.Range("C9").Value = "Report Covering the Week  3/6/2023  Through  3/12/2023"

# --- Crime-complaint figures (rows 15-27) ---------------------------

# Cells that become literal text ("0" / "***.*")
Set-TextValue "F15" "0"
Set-TextValue "C16" "0"
Set-TextValue "D16" "0"
Set-TextValue "E16" "***.*"
Set-TextValue "D23" "0"
Set-TextValue "E23" "***.*"
Set-TextValue "F26" "0"
Set-TextValue "F27" "0"

# Cells that become real numbers again (were literal text before)
Set-NumValue "D18" 1 "#,##0"
Set-NumValue "E18" 0 "#,##0.0;""-""#,##0.0"
Set-NumValue "D20" 1 "#,##0"
Set-NumValue "E20" 0 "#,##0.0;""-""#,##0.0"

# Plain numeric updates (style unchanged)
Set-PlainValue "H15" -100
Set-PlainValue "F16" 3
Set-PlainValue "G16" 2
Set-PlainValue "H16" 50
Set-PlainValue "N16" -85.185185185185
Set-PlainValue "F17" 5
Set-PlainValue "G17" 4
Set-PlainValue "H17" 25
Set-PlainValue "I17" 18
Set-PlainValue "J17" 15
Set-PlainValue "K17" 20
Set-PlainValue "L17" 28.571428571428
Set-PlainValue "M17" 20
Set-PlainValue "N17" -47.058823529411
Set-PlainValue "F18" 3
Set-PlainValue "G18" 3
Set-PlainValue "I18" 8
Set-PlainValue "J18" 5
Set-PlainValue "K18" 60
Set-PlainValue "L18" 14.285714285714
Set-PlainValue "M18" -27.272727272727
Set-PlainValue "N18" -91.919191919191
Set-PlainValue "C19" 1
Set-PlainValue "E19" -50
Set-PlainValue "F19" 10
Set-PlainValue "G19" 8
Set-PlainValue "H19" 25
Set-PlainValue "J19" 27
Set-PlainValue "K19" 11.111111111111
Set-PlainValue "L19" 30.434782608695
Set-PlainValue "M19" 57.894736842105
Set-PlainValue "N19" 11.111111111111
Set-PlainValue "F20" 5
Set-PlainValue "G20" 2
Set-PlainValue "H20" 150
Set-PlainValue "I20" 7
Set-PlainValue "J20" 9
Set-PlainValue "K20" -22.222222222222
Set-PlainValue "L20" 133.333333333333
Set-PlainValue "M20" 75
Set-PlainValue "N20" -82.5
Set-PlainValue "C21" 3
Set-PlainValue "E21" -40
Set-PlainValue "F21" 26
Set-PlainValue "H21" 23.809523809523
Set-PlainValue "I21" 72
Set-PlainValue "J21" 66
Set-PlainValue "K21" 9.090909090909
Set-PlainValue "L21" 38.461538461538
Set-PlainValue "M21" 12.5
Set-PlainValue "N21" -72.413793103448
Set-PlainValue "C24" 7
Set-PlainValue "E24" 40
Set-PlainValue "F24" 34
Set-PlainValue "G24" 29
Set-PlainValue "H24" 17.241379310344
Set-PlainValue "I24" 102
Set-PlainValue "J24" 77
Set-PlainValue "K24" 32.467532467532
Set-PlainValue "L24" 56.923076923076
Set-PlainValue "M24" 78.947368421052
Set-PlainValue "C25" 1
Set-PlainValue "D25" 3
Set-PlainValue "E25" -66.666666666666
Set-PlainValue "F25" 14
Set-PlainValue "H25" 40
Set-PlainValue "I25" 34
Set-PlainValue "J25" 29
Set-PlainValue "K25" 17.241379310344
Set-PlainValue "L25" 13.333333333333
Set-PlainValue "M25" -40.350877192982
Set-PlainValue "H26" -100
Set-PlainValue "G27" 3
Set-PlainValue "H27" -100
Set-PlainValue "J27" 6
Set-PlainValue "K27" -50

